# Split the run that currently reads
#   "...prompt the Goo Interpreter to translate "
# into several runs spelling "...to " + "t" + "co" + "ranslate" + " ",
# wrapping the (now mis-spelled) "tcoranslate" in spellcheck proofErr
# markers - mirrors Word's own behaviour of bracketing a flagged word
# with <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>
# after it gets edited into something the spellchecker doesn't know.

$d = $word.ActiveDocument

# Locate the paragraph that holds the sentence we need to touch.
$needle = "Adding a script to an object does nothing. In order to use it, " +
          "first you need to prompt the Goo Interpreter to translate "
$rng = $d.Content
$found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target sentence - document may have changed."
}

$para = $rng.Paragraphs(1)
$paraRange = $para.Range

$runPr = "<w:rPr>" +
         "<w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/>" +
         "<w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr>"

# Everything up to - and including - the trailing space of "to " stays as
# one run; "translate " is then rebuilt as t / co / ranslate / " " with
# proofErr markers bracketing the (mis-)spelled word, exactly as Word
# would re-flow the paragraph after that kind of in-place typing edit.
$newRuns =
    "<w:r>$runPr<w:t xml:space=`"preserve`">Adding a script to an object " +
    "does nothing. In order to use it, first you need to prompt the Goo " +
    "Interpreter to </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r>$runPr<w:t>t</w:t></w:r>" +
    "<w:r>$runPr<w:t>co</w:t></w:r>" +
    "<w:r>$runPr<w:t>ranslate</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r>$runPr<w:t xml:space=`"preserve`"> </w:t></w:r>"

$oldRun = "<w:r>$runPr<w:t xml:space=`"preserve`">Adding a script to an " +
          "object does nothing. In order to use it, first you need to " +
          "prompt the Goo Interpreter to translate </w:t></w:r>"

# Rebuild the paragraph's OOXML: swap only the opening run for the
# freshly split run sequence, leave every other run in the paragraph
# (the "your", " script into events ...", " method like this:" runs)
# completely untouched.
$paraPrefix = "<w:pPr><w:pStyle w:val=`"NormalWeb`"/>" +
              "<w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"4`"/></w:numPr>" +
              "<w:spacing w:before=`"0`" w:beforeAutospacing=`"0`" " +
              "w:after=`"160`" w:afterAutospacing=`"0`"/>" +
              "<w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" " +
              "w:cs=`"Arial`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/>" +
              "</w:rPr></w:pPr>"

$tailRuns =
    "<w:r w:rsidR=`"00B00106`">$runPr<w:t>your</w:t></w:r>" +
    "<w:r>$runPr<w:t xml:space=`"preserve`"> script into events with " +
    [char]0x201C + "translate()" + [char]0x201D + "</w:t></w:r>" +
    "<w:r w:rsidR=`"00164D18`">$runPr<w:t xml:space=`"preserve`"> method " +
    "like this:</w:t></w:r>"

$newParaInner = $newRuns + $tailRuns

$xml = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
       'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
       '<w:p w14:paraId="50C2B60E" w14:textId="77777777" w:rsidR="001A0FEC" ' +
       'w:rsidRDefault="001A0FEC" w:rsidP="003E2F42">' + $paraPrefix +
       $newParaInner + '</w:p></w:document>'

$paraRange.InsertXML($xml)

Write-Output "Paragraph rebuilt: $($paraRange.Text)"
